# matcg mora vs pera
# Update the point table for the match result between "Mora A" (row 2)
# and "Pera" (row 4): both teams played one more game, which was drawn
# 1-1 (each team scored one more goal and conceded one more goal), and
# Mora A picked up a green card during the match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Mora A
$ws.Range("B2").Value = 4   # Played
$ws.Range("D2").Value = 2   # Drawn
$ws.Range("F2").Value = 6   # Goal Scored
$ws.Range("G2").Value = 1   # Goal Concedered
$ws.Range("L2").Value = 2   # Green card points

# Row 4 - Pera
$ws.Range("B4").Value = 3   # Played
$ws.Range("D4").Value = 2   # Drawn
$ws.Range("F4").Value = 3   # Goal Scored
$ws.Range("G4").Value = 1   # Goal Concedered

$ws.Range("D21").Select()
